# Auto-generated Excel COM-interop script that updates cached Sheet values
# (currentAveragePrice / LevePrice / LeveProfit columns) per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 471.56668
$ws.Range("J17").Value = 439.7037
$ws.Range("L17").Value = 1319.1111
$ws.Range("N17").Value = -1655.1111

# Row 112
$ws.Range("H112").Value = 4422.4116
$ws.Range("J112").Value = 4680.0625
$ws.Range("L112").Value = 14040.1875
$ws.Range("N112").Value = -16256.1875

# Row 125
$ws.Range("H125").Value = 784.6667
$ws.Range("I125").Value = 818
$ws.Range("J125").Value = 718
$ws.Range("K125").Value = 7362
$ws.Range("L125").Value = 6462
$ws.Range("M125").Value = -4902
$ws.Range("N125").Value = -11382

# Row 137
$ws.Range("H137").Value = 904.2692
$ws.Range("I137").Value = 775.55
$ws.Range("J137").Value = 1333.3334
$ws.Range("K137").Value = 2326.65
$ws.Range("L137").Value = 4000.0002
$ws.Range("M137").Value = 223.3500000000004
$ws.Range("N137").Value = -9100.0002

# Row 138
$ws.Range("H138").Value = 3337.877
$ws.Range("J138").Value = 4362.7617
$ws.Range("L138").Value = 13088.2851
$ws.Range("N138").Value = -23368.2851

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 2058.1428
$ws.Range("I26").Value = 1681.4
$ws.Range("K26").Value = 1681.4
$ws.Range("M26").Value = -1351.4

# Row 32
$ws.Range("H32").Value = 327207.25
$ws.Range("I32").Value = 1768.4147
$ws.Range("J32").Value = 4774871.5
$ws.Range("K32").Value = 1768.4147
$ws.Range("L32").Value = 4774871.5
$ws.Range("M32").Value = -1481.4147
$ws.Range("N32").Value = -4775445.5

# Row 41
$ws.Range("H41").Value = 2167.6667
$ws.Range("I41").Value = 2167.6667
$ws.Range("K41").Value = 2167.6667
$ws.Range("M41").Value = -1753.6667

# Row 61
$ws.Range("H61").Value = 3954144.2
$ws.Range("I61").Value = 6994322
$ws.Range("J61").Value = 1912.7
$ws.Range("K61").Value = 6994322
$ws.Range("L61").Value = 1912.7
$ws.Range("M61").Value = -6994110
$ws.Range("N61").Value = -2336.7

# Row 74
$ws.Range("H74").Value = 1235.0488
$ws.Range("I74").Value = 1082.9395
$ws.Range("J74").Value = 1862.5
$ws.Range("K74").Value = 1082.9395
$ws.Range("L74").Value = 1862.5
$ws.Range("M74").Value = -208.9395
$ws.Range("N74").Value = -3610.5

# Row 77
$ws.Range("H77").Value = 1235.0488
$ws.Range("I77").Value = 1082.9395
$ws.Range("J77").Value = 1862.5
$ws.Range("K77").Value = 5414.6975
$ws.Range("L77").Value = 9312.5
$ws.Range("M77").Value = -1046.6975
$ws.Range("N77").Value = -18048.5

# Row 132
$ws.Range("H132").Value = 55612976
$ws.Range("I132").Value = 142858160
$ws.Range("J132").Value = 93308.37
$ws.Range("K132").Value = 428574480
$ws.Range("L132").Value = 279925.11
$ws.Range("M132").Value = -428571950
$ws.Range("N132").Value = -284985.11

# Row 136
$ws.Range("H136").Value = 3954144.2
$ws.Range("I136").Value = 6994322
$ws.Range("J136").Value = 1912.7
$ws.Range("K136").Value = 20982966
$ws.Range("L136").Value = 5738.1
$ws.Range("M136").Value = -20980416
$ws.Range("N136").Value = -10838.1

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 6731.892
$ws.Range("I134").Value = 2924.348
$ws.Range("J134").Value = 12987.143
$ws.Range("K134").Value = 8773.044
$ws.Range("L134").Value = 38961.429
$ws.Range("M134").Value = -6238.044
$ws.Range("N134").Value = -44031.429

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1226.96
$ws.Range("I16").Value = 1125.4
$ws.Range("J16").Value = 1379.3
$ws.Range("K16").Value = 1125.4
$ws.Range("L16").Value = 1379.3
$ws.Range("M16").Value = -838.4000000000001
$ws.Range("N16").Value = -1953.3

# Row 31
$ws.Range("H31").Value = 3248291
$ws.Range("I31").Value = 4311565
$ws.Range("J31").Value = 2506.7368
$ws.Range("K31").Value = 4311565
$ws.Range("L31").Value = 2506.7368
$ws.Range("M31").Value = -4311270
$ws.Range("N31").Value = -3096.7368

# Row 34
$ws.Range("H34").Value = 3248291
$ws.Range("I34").Value = 4311565
$ws.Range("J34").Value = 2506.7368
$ws.Range("K34").Value = 4311565
$ws.Range("L34").Value = 2506.7368
$ws.Range("M34").Value = -4311363
$ws.Range("N34").Value = -2910.7368

# Row 58
$ws.Range("H58").Value = 1215.5
$ws.Range("I58").Value = 1143.1765
$ws.Range("J58").Value = 1327.2727
$ws.Range("K58").Value = 1143.1765
$ws.Range("L58").Value = 1327.2727
$ws.Range("M58").Value = -940.1765
$ws.Range("N58").Value = -1733.2727

# Row 80
$ws.Range("H80").Value = 16666.666
$ws.Range("J80").Value = 16666.666
$ws.Range("L80").Value = 16666.666
$ws.Range("N80").Value = -18912.666

# Row 83
$ws.Range("H83").Value = 16666.666
$ws.Range("J83").Value = 16666.666
$ws.Range("L83").Value = 49999.99800000001
$ws.Range("N83").Value = -61231.99800000001

# Row 113
$ws.Range("H113").Value = 1226.96
$ws.Range("I113").Value = 1125.4
$ws.Range("J113").Value = 1379.3
$ws.Range("K113").Value = 1125.4
$ws.Range("L113").Value = 1379.3
$ws.Range("M113").Value = 1044.6
$ws.Range("N113").Value = -5719.3

# Row 132
$ws.Range("H132").Value = 28652.447
$ws.Range("I132").Value = 1558.7307
$ws.Range("J132").Value = 87355.5
$ws.Range("K132").Value = 4676.1921
$ws.Range("L132").Value = 262066.5
$ws.Range("M132").Value = -2146.1921
$ws.Range("N132").Value = -267126.5

# Row 134
$ws.Range("H134").Value = 1770.3823
$ws.Range("I134").Value = 1324.2963
$ws.Range("J134").Value = 3491
$ws.Range("K134").Value = 3972.8889
$ws.Range("L134").Value = 10473
$ws.Range("M134").Value = -1437.8889
$ws.Range("N134").Value = -15543

# Row 136
$ws.Range("H136").Value = 1215.5
$ws.Range("I136").Value = 1143.1765
$ws.Range("J136").Value = 1327.2727
$ws.Range("K136").Value = 3429.5295
$ws.Range("L136").Value = 3981.8181
$ws.Range("M136").Value = -879.5295000000001
$ws.Range("N136").Value = -9081.8181

$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value = 875375
$ws.Range("J33").Value = 875375
$ws.Range("L33").Value = 875375
$ws.Range("N33").Value = -875879

# Row 102
$ws.Range("H102").Value = 1381.7916
$ws.Range("I102").Value = 1046.6666
$ws.Range("J102").Value = 1716.9166
$ws.Range("K102").Value = 1046.6666
$ws.Range("L102").Value = 1716.9166
$ws.Range("M102").Value = 575.3334
$ws.Range("N102").Value = -4960.9166

# Row 129
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

# Row 132
$ws.Range("H132").Value = 649077.4
$ws.Range("I132").Value = 112411.445
$ws.Range("J132").Value = 1252826.5
$ws.Range("K132").Value = 337234.335
$ws.Range("L132").Value = 3758479.5
$ws.Range("M132").Value = -334704.335
$ws.Range("N132").Value = -3763539.5

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 15886.723
$ws.Range("I132").Value = 22608.918
$ws.Range("J132").Value = 1565.5217
$ws.Range("K132").Value = 67826.754
$ws.Range("L132").Value = 4696.5651
$ws.Range("M132").Value = -65296.754
$ws.Range("N132").Value = -9756.5651

# Row 136
$ws.Range("H136").Value = 5071.122
$ws.Range("I136").Value = 5857.52
$ws.Range("J136").Value = 3842.375
$ws.Range("K136").Value = 17572.56
$ws.Range("L136").Value = 11527.125
$ws.Range("M136").Value = -15022.56
$ws.Range("N136").Value = -16627.125

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 2832.5
$ws.Range("I2").Value = 2351
$ws.Range("K2").Value = 2351
$ws.Range("M2").Value = -2239

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""

# Row 130
$ws.Range("H130").Value = 30333
$ws.Range("J130").Value = 30333
$ws.Range("L130").Value = 30333
$ws.Range("N130").Value = -40373

# Row 132
$ws.Range("H132").Value = 67159250
$ws.Range("I132").Value = 102319700
$ws.Range("J132").Value = 2698445.2
$ws.Range("K132").Value = 306959100
$ws.Range("L132").Value = 8095335.600000001
$ws.Range("M132").Value = -306956570
$ws.Range("N132").Value = -8100395.600000001

# Row 136
$ws.Range("H136").Value = 25022.547
$ws.Range("I136").Value = 44186.39
$ws.Range("J136").Value = 1824.2106
$ws.Range("K136").Value = 132559.17
$ws.Range("L136").Value = 5472.6318
$ws.Range("M136").Value = -130009.17
$ws.Range("N136").Value = -10572.6318
